{"js": "// Update the Step 3.1 observation paragraph:\n//  - \"only the number 20 was sent to the server.\" becomes a multi-run\n//    sentence: \"only the number 21 was sent to the server; indicating\n//    that the loop finished before any threads sent their corresponding\n//    task count.\" (split across several runs, matching the author's\n//    incremental edit pattern)\n//  - a <w:lastRenderedPageBreak/> marker is added at the start of the\n//    paragraph's first run (Word stamps this in when the paragraph\n//    happens to fall at a page boundary on save).\n//\n// We locate the target paragraph by its distinctive text, then replace\n// its contents in one shot via OOXML (Range.insertOoxml), which lets us\n// control the exact run layout instead of relying on plain text-replace\n// (which would collapse everything back into a single run).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst needle = \"We observed there was no problem connecting to the server\";\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(needle) !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the Step 3.1 observation paragraph.\");\n}\n\n// Grab the paragraph's own <w:p ...> opening tag so we can preserve its\n// paraId/rsid attributes on the replacement paragraph.\nconst existing = target.getOoxml();\nawait context.sync();\nconst tagMatch = existing.value.match(/<w:p\\b[^>]*>/);\nconst pOpenTag = tagMatch ? tagMatch[0] : \"<w:p>\";\n\nconst ooxml = `<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          ${pOpenTag}\n            <w:r>\n              <w:lastRenderedPageBreak/>\n              <w:t xml:space=\"preserve\">We observed there was no problem connecting to the server from the client, and no difficulty sending and receiving data back. However, </w:t>\n            </w:r>\n            <w:r>\n              <w:t>only the number 2</w:t>\n            </w:r>\n            <w:r>\n              <w:t>1</w:t>\n            </w:r>\n            <w:r>\n              <w:t xml:space=\"preserve\"> was sent to the server</w:t>\n            </w:r>\n            <w:r>\n              <w:t>; indicating that the loop finished before any threads sent their corresponding task count</w:t>\n            </w:r>\n            <w:r>\n              <w:t>.</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntarget.getRange().insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Update the Step 3.1 observation paragraph:\n#  - \"only the number 20 was sent to the server.\" becomes a multi-run\n#    sentence: \"only the number 21 was sent to the server; indicating\n#    that the loop finished before any threads sent their corresponding\n#    task count.\" (split across several runs, matching the author's\n#    incremental edit pattern)\n#  - a <w:lastRenderedPageBreak/> marker is added at the start of the\n#    paragraph's first run (Word stamps this in when the paragraph\n#    happens to fall at a page boundary on save).\n#\n# We locate the target paragraph by its distinctive text, then replace\n# its contents in one shot via OOXML (Range.InsertXML), which lets us\n# control the exact run layout instead of relying on plain text-replace\n# (which would collapse everything back into a single run).\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*We observed there was no problem connecting to the server*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not find the Step 3.1 observation paragraph.\"\n}\n\n# Grab the paragraph's own <w:p ...> opening tag so we can preserve its\n# paraId/rsid attributes on the replacement paragraph.\n$existingXml = $target.Range.WordOpenXML\n$pOpenTag = \"<w:p>\"\nif ($existingXml -match '<w:p\\b[^>]*>') {\n    $pOpenTag = $matches[0]\n}\n\n$ooxml = @\"\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          $pOpenTag\n            <w:r>\n              <w:lastRenderedPageBreak/>\n              <w:t xml:space=\"preserve\">We observed there was no problem connecting to the server from the client, and no difficulty sending and receiving data back. However, </w:t>\n            </w:r>\n            <w:r>\n              <w:t>only the number 2</w:t>\n            </w:r>\n            <w:r>\n              <w:t>1</w:t>\n            </w:r>\n            <w:r>\n              <w:t xml:space=\"preserve\"> was sent to the server</w:t>\n            </w:r>\n            <w:r>\n              <w:t>; indicating that the loop finished before any threads sent their corresponding task count</w:t>\n            </w:r>\n            <w:r>\n              <w:t>.</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$target.Range.InsertXML($ooxml)\n"}
